# Add new weekly ranking sheet: magapoke_2026-01-07
$wb = $excel.ActiveWorkbook

$sheetName = "magapoke_2026-01-07"
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)

# Insert the new sheet right after the current last sheet (end of tab strip)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = $sheetName

# Header row
$ws.Cells.Item(1, 1).Value = "rank"
$ws.Cells.Item(1, 2).Value = "title"

# Reuse the header formatting (bold, centered, thin border) from the
# previous week's sheet so the new header lands on the same shared style.
$lastSheet.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the page margins used by every other weekly sheet (0.75in / 1in / 0.5in).
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

$titles = @(
    'ブルーロック',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    'みいちゃんと山田さん',
    'ガチアクタ',
    '東京卍リベンジャーズ',
    'ベイビーステップ',
    'ドラハチ',
    '島耕作',
    'FAIRY TAIL 100 YEARS QUEST',
    '薫る花は凛と咲く',
    'イレギュラーズ',
    '黄昏町プリズナーズ',
    '魔女と傭兵',
    '十字架のろくにん',
    'ハードワーカー中田',
    '君が僕らを悪魔と呼んだ頃',
    'WIND BREAKER',
    'ひゃくえむ。',
    'ギルティサークル',
    'GALAXIAS',
    'K-9~警視庁公安部公安第9課異能対策係~',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    'アルキメデスの大戦',
    '蒼く染めろ',
    '幼馴染とはラブコメにならない',
    '愛妻の裏アカ',
    'さわらないで小手指くん',
    'せいぶつ部の田辺くん',
    'グラぱらっ！',
    'ハナバス　苔石花江のバスケ論',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    '屋根の下のアルテミス',
    '食糧人類-Starving Anonymous-',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '南海トラフ巨大地震',
    'デッドアカウント',
    '異世界ウォーキング',
    'アオバノバスケ',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    '春くらり',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    'となりの黒川さん',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    'いじめるヤバイ奴',
    '普通の本はありません！',
    '五輪の女神さま 〜なでしこ寮のメダルごはん〜',
    'おやすみ ふみさん',
    'なれの果ての僕ら',
    '黒猫と魔女の教室',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    '降り積もれ孤独な死よ',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    '阿武ノーマル',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '金田一少年の事件簿外伝 犯人たちの事件簿',
    'デスティニーラバーズ',
    'MYS',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    'ともだちづくり',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    'ヒロインは絶望しました。',
    '我間乱 ―修羅―',
    '人間消失',
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
    '君が監督！',
    'ジュミドロ',
    'ストーカー行為がバレて人生終了男',
    'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～',
    '魁の花巫女',
    '中華一番！極',
    'イジらないで、長瀞さん',
    '花園さんちのふたごちゃん',
    '日本語が話せないロシア人美少女転入生が頼れるのは、多言語マスターの俺1人',
    'ナキナギ',
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
    '東京ネオンスキャンダル',
    'インフェクション',
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
    '可愛いだけじゃない式守さん',
    '放課後、ぼくは君になる',
    'Fate/Grand Order -Epic of Remnant- 英霊剣豪七番勝負',
    '〈小市民〉 春期限定いちごタルト事件',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    'SHAMAN KING THE SUPER STAR',
    '微妙に優しいいじめっ子',
    'それがメイドのカンナです',
    'ウイニング パス',
    'DAYS外伝',
    '鳴るさんだぁ',
    '東京デスレース',
    '勇者と呼ばれた後に　―そして無双男は家族を創る―'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

$ws.Range("A1").Select()

